$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("Y2").Value = 7628.97
$ws.Range("AG2").Value = 205774.05

# Row 3 - Bibi Cell Vieiralves
$ws.Range("Y3").Value = 3791
$ws.Range("AG3").Value = 104594.01

# Row 4 - Bibi Cell Manauara
$ws.Range("X4").Value = 2636
$ws.Range("Y4").Value = 2563
$ws.Range("AG4").Value = 71424.89999999999

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("Y5").Value = 1047
$ws.Range("AG5").Value = 61543.79

# Row 6 - total
$ws.Range("X6").Value = 14800.32
$ws.Range("Y6").Value = 15029.97
$ws.Range("AG6").Value = 443336.75
